$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "200ohm ferrite bead"
$ws.Range("F5").Value = "200ohm"
$ws.Range("G5:J5").ClearContents()
